$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "NegativeLoginTest"

$ws.Range("A1").Value = "Student-5"
$ws.Range("B1").Value = "S1234"
$ws.Range("A2").Value = "  "
$ws.Range("B2").Value = "S12345"
$ws.Range("A3").Value = "S12345"
$ws.Range("B3").Value = "Student_12345"
$ws.Range("A4").Value = "Student_5"
$ws.Range("B4").Value = " "
$ws.Range("A5").Value = "  "
$ws.Range("B5").Value = "  "
$ws.Range("A6").Value = "S12345"
$ws.Range("B6").Value = "  "

$ws.Columns.Item(1).ColumnWidth = 11.5
$ws.Columns.Item(2).ColumnWidth = 15

$null = $ws.Range("B6").Select()
